$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 144.783305
$ws.Range("H2").Value = 434.349915
$ws.Range("I2").Value = 0.2430046335191003
$ws.Range("J2").Value = 0.251012682214973
$ws.Range("O2").Value = 0.5489460219708985
$ws.Range("P2").Value = 0.6460858926464538
$ws.Range("Q2").Value = 57.21570777540833
$ws.Range("R2").Value = 514.941369978675
$ws.Range("S2").Value = 0.1333964268908062
$ws.Range("T2").Value = 0.1621757528544414
$ws.Range("G3").Value = 144.783305
$ws.Range("H3").Value = 434.349915
$ws.Range("I3").Value = 0.2430046335191003
$ws.Range("J3").Value = 0.251012682214973
$ws.Range("M3").Value = 0.32471
$ws.Range("N3").Value = 0.64942
$ws.Range("O3").Value = 0.4510539780291016
$ws.Range("P3").Value = 0.3539141073535463
$ws.Range("Q3").Value = 47.01258696655
$ws.Range("R3").Value = 282.0755217993
$ws.Range("S3").Value = 0.1096082066282942
$ws.Range("T3").Value = 0.08883692936053154
$ws.Range("G4").Value = 82.24887099999999
$ws.Range("I4").Value = 0.1380466950572427
$ws.Range("J4").Value = 0.1425959278859072
$ws.Range("O4").Value = 0.5489460219708985
$ws.Range("P4").Value = 0.6460858926464538
$ws.Range("S4").Value = 0.07578018409790306
$ws.Range("T4").Value = 0.0921292173559157
$ws.Range("G5").Value = 82.24887099999999
$ws.Range("I5").Value = 0.1380466950572427
$ws.Range("J5").Value = 0.1425959278859072
$ws.Range("M5").Value = 0.32471
$ws.Range("N5").Value = 0.64942
$ws.Range("O5").Value = 0.4510539780291016
$ws.Range("P5").Value = 0.3539141073535463
$ws.Range("Q5").Value = 26.70703090241
$ws.Range("R5").Value = 160.24218541446
$ws.Range("S5").Value = 0.06226651095933963
$ws.Range("T5").Value = 0.0504667105299915
$ws.Range("G6").Value = 163.8590903333333
$ws.Range("H6").Value = 491.577271
$ws.Range("I6").Value = 0.2750214756820535
$ws.Range("J6").Value = 0.284084617144743
$ws.Range("O6").Value = 0.5489460219708985
$ws.Range("P6").Value = 0.6460858926464538
$ws.Range("Q6").Value = 64.75410841641056
$ws.Range("R6").Value = 582.786975747695
$ws.Range("S6").Value = 0.1509719450322295
$ws.Range("T6").Value = 0.1835430634550873
$ws.Range("G7").Value = 163.8590903333333
$ws.Range("H7").Value = 491.577271
$ws.Range("I7").Value = 0.2750214756820535
$ws.Range("J7").Value = 0.284084617144743
$ws.Range("M7").Value = 0.32471
$ws.Range("N7").Value = 0.64942
$ws.Range("O7").Value = 0.4510539780291016
$ws.Range("P7").Value = 0.3539141073535463
$ws.Range("Q7").Value = 53.20668522213667
$ws.Range("R7").Value = 319.24011133282
$ws.Range("S7").Value = 0.124049530649824
$ws.Range("T7").Value = 0.1005415536896557
$ws.Range("G8").Value = 57.0238095
$ws.Range("H8").Value = 114.047619
$ws.Range("I8").Value = 0.09570889357312636
$ws.Range("J8").Value = 0.06590860906562239
$ws.Range("O8").Value = 0.5489460219708985
$ws.Range("P8").Value = 0.6460858926464538
$ws.Range("Q8").Value = 22.5347640778925
$ws.Range("R8").Value = 135.208584467355
$ws.Range("S8").Value = 0.0525390163942038
$ws.Range("T8").Value = 0.04258262252124879
$ws.Range("G9").Value = 57.0238095
$ws.Range("H9").Value = 114.047619
$ws.Range("I9").Value = 0.09570889357312636
$ws.Range("J9").Value = 0.06590860906562239
$ws.Range("M9").Value = 0.32471
$ws.Range("N9").Value = 0.64942
$ws.Range("O9").Value = 0.4510539780291016
$ws.Range("P9").Value = 0.3539141073535463
$ws.Range("Q9").Value = 18.516201182745
$ws.Range("R9").Value = 74.06480473098
$ws.Range("S9").Value = 0.04316987717892255
$ws.Range("T9").Value = 0.02332598654437359
$ws.Range("G10").Value = 147.8896333333333
$ws.Range("H10").Value = 443.6689
$ws.Range("I10").Value = 0.2482183021684772
$ws.Range("J10").Value = 0.2563981636887546
$ws.Range("O10").Value = 0.5489460219708985
$ws.Range("P10").Value = 0.6460858926464538
$ws.Range("Q10").Value = 58.44327178338889
$ws.Range("R10").Value = 525.9894460505001
$ws.Range("S10").Value = 0.136258449555756
$ws.Range("T10").Value = 0.1656552364597606
$ws.Range("G11").Value = 147.8896333333333
$ws.Range("H11").Value = 443.6689
$ws.Range("I11").Value = 0.2482183021684772
$ws.Range("J11").Value = 0.2563981636887546
$ws.Range("M11").Value = 0.32471
$ws.Range("N11").Value = 0.64942
$ws.Range("O11").Value = 0.4510539780291016
$ws.Range("P11").Value = 0.3539141073535463
$ws.Range("Q11").Value = 48.02124283966667
$ws.Range("R11").Value = 288.127457038
$ws.Range("S11").Value = 0.1119598526127212
$ws.Range("T11").Value = 0.09074292722899402
